$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price text is a plain decimal number (single dot, e.g. "245.70")
# need to be explicitly formatted as Text first, otherwise Excel's normal
# value-entry coercion turns them into numbers, which silently drops trailing
# zeros (245.70 -> 245.7) or renders tiny values in scientific notation
# (0.000007939 -> 7.939E-06). The source column stores these as literal text,
# so force the same behaviour here (values with two dots, e.g. "30.347.20",
# are never parsed as numbers by Excel and need no special handling).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range('D2').Value = '30.347.20'
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').Value = '1.938.34'
$ws.Range('E3').Value = '  +0.07%  '
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  +0.34%  '
$ws.Range('D5').Value = '0.7683'
$ws.Range('E5').Value = '  +6.21%  '
$ws.Range('D6').Value = '245.70'
$ws.Range('E6').Value = '  -2.48%  '
$ws.Range('D7').Value = '1.003'
$ws.Range('E7').Value = '  +0.25%  '
$ws.Range('E8').Value = '  -3.34%  '
$ws.Range('D9').Value = '27.68'
$ws.Range('E9').Value = '  -1.13%  '
$ws.Range('D10').Value = '0.07019'
$ws.Range('E10').Value = '  -3.40%  '
$ws.Range('D11').Value = '0.7811'
$ws.Range('E11').Value = '  -3.50%  '
$ws.Range('D12').Value = '0.08007'
$ws.Range('E12').Value = '  -1.21%  '
$ws.Range('D13').Value = '1.933.34'
$ws.Range('E13').Value = '  -0.12%  '
$ws.Range('D14').Value = '5.352'
$ws.Range('E14').Value = '  -2.47%  '
$ws.Range('D15').Value = '94.56'
$ws.Range('E15').Value = '  -0.44%  '
$ws.Range('D16').Value = '14.47'
$ws.Range('E16').Value = '  -4.32%  '
$ws.Range('D17').Value = '30.353.23'
$ws.Range('E17').Value = '  +0.04%  '
$ws.Range('D18').Value = '255.20'
$ws.Range('E18').Value = '  +0.74%  '
$ws.Range('D19').Value = '0.000007939'
$ws.Range('E19').Value = '  -3.71%  '
$ws.Range('D20').Value = '5.775'
$ws.Range('E20').Value = '  -1.11%  '
$ws.Range('D21').Value = '2.189.14'
$ws.Range('E21').Value = '  -0.07%  '
$ws.Range('E22').Value = '  +0.24%  '
$ws.Range('D23').Value = '1.004'
$ws.Range('E23').Value = '  +0.35%  '
$ws.Range('D24').Value = '6.695'
$ws.Range('E24').Value = '  -3.86%  '
$ws.Range('D25').Value = '9.520'
$ws.Range('E25').Value = '  -2.55%  '
$ws.Range('D26').Value = '164.72'
$ws.Range('E26').Value = '  -0.78%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '19.06'
$ws.Range('E27').Value = '  -1.56%  '
$ws.Range('B28').Value = 'Stellar'
$ws.Range('C28').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D28').Value = '0.1337'
$ws.Range('E28').Value = '  +2.80%  '
$ws.Range('D29').Value = '2.271'
$ws.Range('E29').Value = '  -3.40%  '
$ws.Range('E30').Value = '  +0.89%  '
$ws.Range('E31').Value = '  -1.95%  '
$ws.Range('D32').Value = '4.400'
$ws.Range('E32').Value = '  -1.01%  '
$ws.Range('D33').Value = '4.119'
$ws.Range('E33').Value = '  -2.27%  '
$ws.Range('D34').Value = '0.05154'
$ws.Range('E34').Value = '  -1.92%  '
$ws.Range('D35').Value = '1.276'
$ws.Range('E35').Value = '  +0.51%  '
$ws.Range('D36').Value = '0.7473'
$ws.Range('E36').Value = '  -0.59%  '
$ws.Range('D37').Value = '2.788'
$ws.Range('E37').Value = '  +0.64%  '
$ws.Range('D38').Value = '0.01952'
$ws.Range('E38').Value = '  -1.11%  '
$ws.Range('D39').Value = '2.810'
$ws.Range('E39').Value = '  +0.19%  '
$ws.Range('D40').Value = '78.74'
$ws.Range('E40').Value = '  -0.99%  '
$ws.Range('D41').Value = '6.416'
$ws.Range('E41').Value = '  -0.56%  '
$ws.Range('D42').Value = '0.4487'
$ws.Range('E42').Value = '  -1.46%  '
$ws.Range('D43').Value = '1.969'
$ws.Range('E44').Value = '  +0.25%  '
$ws.Range('D45').Value = '0.8339'
$ws.Range('E45').Value = '  -1.25%  '
$ws.Range('D46').Value = '101.07'
$ws.Range('E46').Value = '  -0.90%  '
$ws.Range('D47').Value = '9.793'
$ws.Range('E47').Value = '  -0.11%  '
$ws.Range('D48').Value = '7.486'
$ws.Range('E48').Value = '  +0.42%  '
$ws.Range('D49').Value = '978.15'
$ws.Range('E49').Value = '  +10.23%  '
$ws.Range('D50').Value = '37.19'
$ws.Range('E50').Value = '  +1.11%  '
$ws.Range('D51').Value = '0.4146'
$ws.Range('E51').Value = '  -1.28%  '
